$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Backtest_Results")

$ws.Range("C319").Value = 1
$ws.Range("C328").Value = 1
$ws.Range("C329").Value = 1
$ws.Range("C331").Value = 1
$ws.Range("C335").Value = 1
$ws.Range("C336").Value = 1
$ws.Range("C345").Value = 1
$ws.Range("C346").Value = 1
$ws.Range("C350").Value = 1
$ws.Range("C356").Value = 1
$ws.Range("C358").Value = 1
$ws.Range("C362").Value = 1
$ws.Range("C365").Value = 1
$ws.Range("C368").Value = 1
$ws.Range("C372").Value = 1
$ws.Range("C373").Value = 1
$ws.Range("C374").Value = 0
$ws.Range("C375").Value = 0
$ws.Range("C379").Value = 1
$ws.Range("C380").Value = 1
$ws.Range("C383").Value = 1
$ws.Range("C384").Value = 1
$ws.Range("C385").Value = 1
$ws.Range("C388").Value = 1
$ws.Range("C390").Value = 1
$ws.Range("C400").Value = 1
$ws.Range("C410").Value = 1
$ws.Range("C414").Value = 1
$ws.Range("C415").Value = 1
$ws.Range("C418").Value = 1
$ws.Range("C419").Value = 1
$ws.Range("C420").Value = 1
$ws.Range("C423").Value = 1
$ws.Range("C425").Value = 1
$ws.Range("C427").Value = 1
$ws.Range("C431").Value = 1
$ws.Range("C435").Value = 1
$ws.Range("C437").Value = 1
$ws.Range("C440").Value = 1
$ws.Range("C454").Value = 1
$ws.Range("C455").Value = 1
$ws.Range("C456").Value = 1
$ws.Range("C461").Value = 0
$ws.Range("C467").Value = 1
$ws.Range("C468").Value = 1
$ws.Range("C472").Value = 1
$ws.Range("C480").Value = 1
$ws.Range("C481").Value = 1
$ws.Range("C484").Value = 1
$ws.Range("C492").Value = 1
$ws.Range("C494").Value = 1
$ws.Range("C499").Value = 0
$ws.Range("C504").Value = 1
$ws.Range("C511").Value = 0
$ws.Range("C529").Value = 1
$ws.Range("C537").Value = 1
$ws.Range("C540").Value = 1
$ws.Range("C541").Value = 0
$ws.Range("C544").Value = 0
$ws.Range("C548").Value = 0
$ws.Range("C550").Value = 1
$ws.Range("C553").Value = 0
$ws.Range("C556").Value = 1
$ws.Range("C557").Value = 1
$ws.Range("C561").Value = 1
$ws.Range("C564").Value = 1
$ws.Range("C569").Value = 1
$ws.Range("C579").Value = 1
$ws.Range("C583").Value = 0
$ws.Range("C588").Value = 1

$wsm = $wb.Worksheets.Item("Metrics")
$wsm.Range("B2").Value = 0.465195246179966
$wsm.Range("B3").Value = 0.6476510067114094
